$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.315.91"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "2.585.04"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'565.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'142.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "2.591.52"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'0.160"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +11.34%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "3.038.38"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "59.275.88"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "'22.66"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "2.586.44"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").Value = "'4.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "'337.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "'10.38"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'6.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'64.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("E25").Value = "  +5.09%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "'7.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "0.0₃0778"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'161.80"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "'6.08"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'18.95"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "'1.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'0.880"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("D38").Value = "'0.878"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").Value = "'37.55"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "'299.26"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'132.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.60%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "'0.0974"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "'0.597"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'0.0536"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'10.63"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "'19.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "'0.0232"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").Value = "'18.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.62%  "
